$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G2").Value = 17.81824933333333
$ws.Range("H2").Value = 53.454748
$ws.Range("I2").Value = 0.05180179233147032
$ws.Range("J2").Value = 0.05180179233147034
$ws.Range("M2").Value = 0.05237733333333333
$ws.Range("N2").Value = 0.157132
$ws.Range("O2").Value = 0.01249705432598326
$ws.Range("P2").Value = 0.01249705432598326
$ws.Range("Q2").Value = 0.9332723847484444
$ws.Range("R2").Value = 8.399451462736
$ws.Range("S2").Value = 0.0006473698129496876
$ws.Range("T2").Value = 0.0006473698129496878
$ws.Range("G3").Value = 17.81824933333333
$ws.Range("H3").Value = 53.454748
$ws.Range("I3").Value = 0.05180179233147032
$ws.Range("J3").Value = 0.05180179233147034
$ws.Range("O3").Value = 0.1346970137168397
$ws.Range("P3").Value = 0.1346970137168397
$ws.Range("Q3").Value = 10.05909072097422
$ws.Range("R3").Value = 90.531816488768
$ws.Range("S3").Value = 0.006977546732228941
$ws.Range("T3").Value = 0.006977546732228944
$ws.Range("G4").Value = 17.81824933333333
$ws.Range("H4").Value = 53.454748
$ws.Range("I4").Value = 0.05180179233147032
$ws.Range("J4").Value = 0.05180179233147034
$ws.Range("O4").Value = 0.020817156814363
$ws.Range("P4").Value = 0.020817156814363
$ws.Range("Q4").Value = 1.554612557251111
$ws.Range("R4").Value = 13.99151301526
$ws.Range("S4").Value = 0.001078366034229285
$ws.Range("T4").Value = 0.001078366034229285
$ws.Range("G5").Value = 17.81824933333333
$ws.Range("H5").Value = 53.454748
$ws.Range("I5").Value = 0.05180179233147032
$ws.Range("J5").Value = 0.05180179233147034
$ws.Range("M5").Value = 3.48701
$ws.Range("N5").Value = 10.46103
$ws.Range("O5").Value = 0.831988775142814
$ws.Range("P5").Value = 0.831988775142814
$ws.Range("Q5").Value = 62.13241360782667
$ws.Range("R5").Value = 559.19172247044
$ws.Range("S5").Value = 0.04309850975206241
$ws.Range("T5").Value = 0.04309850975206243
$ws.Range("I6").Value = 0.4402211587141748
$ws.Range("J6").Value = 0.4402211587141748
$ws.Range("M6").Value = 0.05237733333333333
$ws.Range("N6").Value = 0.157132
$ws.Range("O6").Value = 0.01249705432598326
$ws.Range("P6").Value = 0.01249705432598326
$ws.Range("Q6").Value = 7.931120374773333
$ws.Range("R6").Value = 71.38008337296
$ws.Range("S6").Value = 0.00550146773589834
$ws.Range("T6").Value = 0.00550146773589834
$ws.Range("I7").Value = 0.4402211587141748
$ws.Range("J7").Value = 0.4402211587141748
$ws.Range("O7").Value = 0.1346970137168397
$ws.Range("P7").Value = 0.1346970137168397
$ws.Range("S7").Value = 0.05929647545376628
$ws.Range("T7").Value = 0.05929647545376629
$ws.Range("I8").Value = 0.4402211587141748
$ws.Range("J8").Value = 0.4402211587141748
$ws.Range("O8").Value = 0.020817156814363
$ws.Range("P8").Value = 0.020817156814363
$ws.Range("S8").Value = 0.00916415289395356
$ws.Range("T8").Value = 0.009164152893953562
$ws.Range("I9").Value = 0.4402211587141748
$ws.Range("J9").Value = 0.4402211587141748
$ws.Range("M9").Value = 3.48701
$ws.Range("N9").Value = 10.46103
$ws.Range("O9").Value = 0.831988775142814
$ws.Range("P9").Value = 0.831988775142814
$ws.Range("Q9").Value = 528.0126783476001
$ws.Range("R9").Value = 4752.1141051284
$ws.Range("S9").Value = 0.3662590626305566
$ws.Range("T9").Value = 0.3662590626305566
$ws.Range("G10").Value = 100.6958183333333
$ws.Range("H10").Value = 302.087455
$ws.Range("I10").Value = 0.29274614875843
$ws.Range("J10").Value = 0.2927461487584301
$ws.Range("M10").Value = 0.05237733333333333
$ws.Range("N10").Value = 0.157132
$ws.Range("O10").Value = 0.01249705432598326
$ws.Range("P10").Value = 0.01249705432598326
$ws.Range("Q10").Value = 5.274178442117778
$ws.Range("R10").Value = 47.46760597906
$ws.Range("S10").Value = 0.003658464524756476
$ws.Range("T10").Value = 0.003658464524756477
$ws.Range("G11").Value = 100.6958183333333
$ws.Range("H11").Value = 302.087455
$ws.Range("I11").Value = 0.29274614875843
$ws.Range("J11").Value = 0.2927461487584301
$ws.Range("O11").Value = 0.1346970137168397
$ws.Range("P11").Value = 0.1346970137168397
$ws.Range("Q11").Value = 56.8466830208089
$ws.Range("R11").Value = 511.6201471872801
$ws.Range("S11").Value = 0.03943203201486625
$ws.Range("T11").Value = 0.03943203201486626
$ws.Range("G12").Value = 100.6958183333333
$ws.Range("H12").Value = 302.087455
$ws.Range("I12").Value = 0.29274614875843
$ws.Range("J12").Value = 0.2927461487584301
$ws.Range("O12").Value = 0.020817156814363
$ws.Range("P12").Value = 0.020817156814363
$ws.Range("Q12").Value = 8.785542323219445
$ws.Range("R12").Value = 79.06988090897501
$ws.Range("S12").Value = 0.006094142485505077
$ws.Range("T12").Value = 0.006094142485505078
$ws.Range("G13").Value = 100.6958183333333
$ws.Range("H13").Value = 302.087455
$ws.Range("I13").Value = 0.29274614875843
$ws.Range("J13").Value = 0.2927461487584301
$ws.Range("M13").Value = 3.48701
$ws.Range("N13").Value = 10.46103
$ws.Range("O13").Value = 0.831988775142814
$ws.Range("P13").Value = 0.831988775142814
$ws.Range("Q13").Value = 351.1273254865167
$ws.Range("R13").Value = 3160.145929378651
$ws.Range("S13").Value = 0.2435615097333022
$ws.Range("T13").Value = 0.2435615097333023
$ws.Range("G14").Value = 74.032918
$ws.Range("H14").Value = 222.098754
$ws.Range("I14").Value = 0.2152309001959248
$ws.Range("J14").Value = 0.2152309001959249
$ws.Range("M14").Value = 0.05237733333333333
$ws.Range("N14").Value = 0.157132
$ws.Range("O14").Value = 0.01249705432598326
$ws.Range("P14").Value = 0.01249705432598326
$ws.Range("Q14").Value = 3.877646823725333
$ws.Range("R14").Value = 34.898821413528
$ws.Range("S14").Value = 0.002689752252378753
$ws.Range("T14").Value = 0.002689752252378754
$ws.Range("G15").Value = 74.032918
$ws.Range("H15").Value = 222.098754
$ws.Range("I15").Value = 0.2152309001959248
$ws.Range("J15").Value = 0.2152309001959249
$ws.Range("O15").Value = 0.1346970137168397
$ws.Range("P15").Value = 0.1346970137168397
$ws.Range("Q15").Value = 41.79444481716266
$ws.Range("R15").Value = 376.150003354464
$ws.Range("S15").Value = 0.02899095951597825
$ws.Range("T15").Value = 0.02899095951597826
$ws.Range("G16").Value = 74.032918
$ws.Range("H16").Value = 222.098754
$ws.Range("I16").Value = 0.2152309001959248
$ws.Range("J16").Value = 0.2152309001959249
$ws.Range("O16").Value = 0.020817156814363
$ws.Range("P16").Value = 0.020817156814363
$ws.Range("Q16").Value = 6.459248707303333
$ws.Range("R16").Value = 58.13323836572999
$ws.Range("S16").Value = 0.00448049540067508
$ws.Range("T16").Value = 0.004480495400675082
$ws.Range("G17").Value = 74.032918
$ws.Range("H17").Value = 222.098754
$ws.Range("I17").Value = 0.2152309001959248
$ws.Range("J17").Value = 0.2152309001959249
$ws.Range("M17").Value = 3.48701
$ws.Range("N17").Value = 10.46103
$ws.Range("O17").Value = 0.831988775142814
$ws.Range("P17").Value = 0.831988775142814
$ws.Range("Q17").Value = 258.15352539518
$ws.Range("R17").Value = 2323.38172855662
$ws.Range("S17").Value = 0.1790696930268927
$ws.Range("T17").Value = 0.1790696930268928
